$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from row 15 (last existing data row) to the two new rows first
$ws.Range("A15:D15").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A17:D17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Set numeric / index cells
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(17, 1).Value = 16

# Populate new shared strings in the exact order they were authored:
# D17, C17, C16, D16, B17, B16
$ws.Cells.Item(17, 4).Value = "Verify the max length allowed at each editText available at Kirana Bazaar registration activity"
$ws.Cells.Item(17, 3).Value = "Registration"
$ws.Cells.Item(16, 3).Value = "Feedback"
$ws.Cells.Item(16, 4).Value = "Verify the max length allowed at feedback editText available at Kirana Bazaar and validate the succesful submission of feedback."
$ws.Cells.Item(17, 2).Value = "Validate_Registration_Maxlength"
$ws.Cells.Item(16, 2).Value = "Validate_Feedback"

# Row heights
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 30

# Update sheet view: scroll + selection
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("D14").Select()
